# Apply scheduled market-data refresh values to the Bahamut profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3465100.2
$ws.Range("I74").Value = 3997115.5
$ws.Range("K74").Value = 3997115.5
$ws.Range("M74").Value = -3996179.5
$ws.Range("H77").Value = 3465100.2
$ws.Range("I77").Value = 3997115.5
$ws.Range("K77").Value = 19985577.5
$ws.Range("M77").Value = -19980897.5
$ws.Range("H80").Value = 432486.5
$ws.Range("I80").Value = 423.41666
$ws.Range("J80").Value = 1296612.6
$ws.Range("K80").Value = 1270.24998
$ws.Range("L80").Value = 3889837.8
$ws.Range("M80").Value = -272.2499800000001
$ws.Range("N80").Value = -3891833.8
$ws.Range("H82").Value = 28576376
$ws.Range("I82").Value = 2635.5
$ws.Range("J82").Value = 66674696
$ws.Range("K82").Value = 7906.5
$ws.Range("L82").Value = 200024088
$ws.Range("M82").Value = -7500.5
$ws.Range("N82").Value = -200024900
$ws.Range("H83").Value = 432486.5
$ws.Range("I83").Value = 423.41666
$ws.Range("J83").Value = 1296612.6
$ws.Range("K83").Value = 3810.74994
$ws.Range("L83").Value = 11669513.4
$ws.Range("M83").Value = 1181.25006
$ws.Range("N83").Value = -11679497.4
$ws.Range("H85").Value = 28576376
$ws.Range("I85").Value = 2635.5
$ws.Range("J85").Value = 66674696
$ws.Range("K85").Value = 7906.5
$ws.Range("L85").Value = 200024088
$ws.Range("M85").Value = -6502.5
$ws.Range("N85").Value = -200026896
$ws.Range("H86").Value = 2676.762
$ws.Range("I86").Value = 1568.5
$ws.Range("J86").Value = 4154.4443
$ws.Range("K86").Value = 1568.5
$ws.Range("L86").Value = 4154.4443
$ws.Range("M86").Value = -445.5
$ws.Range("N86").Value = -6400.4443
$ws.Range("H88").Value = 691687.9399999999
$ws.Range("I88").Value = 2533.8333
$ws.Range("J88").Value = 934918.8
$ws.Range("K88").Value = 2533.8333
$ws.Range("L88").Value = 934918.8
$ws.Range("M88").Value = -2127.8333
$ws.Range("N88").Value = -935730.8
$ws.Range("H89").Value = 2676.762
$ws.Range("I89").Value = 1568.5
$ws.Range("J89").Value = 4154.4443
$ws.Range("K89").Value = 7842.5
$ws.Range("L89").Value = 20772.2215
$ws.Range("M89").Value = -2226.5
$ws.Range("N89").Value = -32004.2215
$ws.Range("H91").Value = 691687.9399999999
$ws.Range("I91").Value = 2533.8333
$ws.Range("J91").Value = 934918.8
$ws.Range("K91").Value = 2533.8333
$ws.Range("L91").Value = 934918.8
$ws.Range("M91").Value = -1129.8333
$ws.Range("N91").Value = -937726.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5483.5713
$ws.Range("I32").Value = 4560.844
$ws.Range("J32").Value = 9258.362999999999
$ws.Range("K32").Value = 4560.844
$ws.Range("L32").Value = 9258.362999999999
$ws.Range("M32").Value = -4273.844
$ws.Range("N32").Value = -9832.362999999999
$ws.Range("H61").Value = 1167.7742
$ws.Range("I61").Value = 1248.05
$ws.Range("J61").Value = 1021.8182
$ws.Range("K61").Value = 1248.05
$ws.Range("L61").Value = 1021.8182
$ws.Range("M61").Value = -1036.05
$ws.Range("N61").Value = -1445.8182
$ws.Range("H74").Value = 1028.091
$ws.Range("I74").Value = 961.8461
$ws.Range("J74").Value = 1123.7778
$ws.Range("K74").Value = 961.8461
$ws.Range("L74").Value = 1123.7778
$ws.Range("M74").Value = -87.84609999999998
$ws.Range("N74").Value = -2871.7778
$ws.Range("H77").Value = 1028.091
$ws.Range("I77").Value = 961.8461
$ws.Range("J77").Value = 1123.7778
$ws.Range("K77").Value = 4809.2305
$ws.Range("L77").Value = 5618.889
$ws.Range("M77").Value = -441.2304999999997
$ws.Range("N77").Value = -14354.889
$ws.Range("H88").Value = 2741.6875
$ws.Range("I88").Value = 2232
$ws.Range("J88").Value = 2973.3635
$ws.Range("K88").Value = 2232
$ws.Range("L88").Value = 2973.3635
$ws.Range("M88").Value = -1826
$ws.Range("N88").Value = -3785.3635
$ws.Range("H91").Value = 2741.6875
$ws.Range("I91").Value = 2232
$ws.Range("J91").Value = 2973.3635
$ws.Range("K91").Value = 2232
$ws.Range("L91").Value = 2973.3635
$ws.Range("M91").Value = -828
$ws.Range("N91").Value = -5781.363499999999
$ws.Range("H97").Value = 524.15
$ws.Range("I97").Value = 487.70587
$ws.Range("J97").Value = 730.6667
$ws.Range("K97").Value = 487.70587
$ws.Range("L97").Value = 730.6667
$ws.Range("M97").Value = 8.294129999999996
$ws.Range("N97").Value = -1722.6667
$ws.Range("H136").Value = 1167.7742
$ws.Range("I136").Value = 1248.05
$ws.Range("J136").Value = 1021.8182
$ws.Range("K136").Value = 3744.15
$ws.Range("L136").Value = 3065.4546
$ws.Range("M136").Value = -1194.15
$ws.Range("N136").Value = -8165.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 58859480
$ws.Range("I20").Value = 60577.7
$ws.Range("J20").Value = 142857920
$ws.Range("K20").Value = 60577.7
$ws.Range("L20").Value = 142857920
$ws.Range("M20").Value = -60330.7
$ws.Range("N20").Value = -142858414
$ws.Range("H86").Value = 2117.0908
$ws.Range("I86").Value = 1787.5555
$ws.Range("J86").Value = 3600
$ws.Range("K86").Value = 1787.5555
$ws.Range("L86").Value = 3600
$ws.Range("M86").Value = -664.5554999999999
$ws.Range("N86").Value = -5846
$ws.Range("H89").Value = 2117.0908
$ws.Range("I89").Value = 1787.5555
$ws.Range("J89").Value = 3600
$ws.Range("K89").Value = 8937.7775
$ws.Range("L89").Value = 18000
$ws.Range("M89").Value = -3321.7775
$ws.Range("N89").Value = -29232
$ws.Range("H105").Value = 5855.316
$ws.Range("I105").Value = 5583.8125
$ws.Range("K105").Value = 5583.8125
$ws.Range("M105").Value = -3836.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1741.1
$ws.Range("I16").Value = 1738.875
$ws.Range("K16").Value = 1738.875
$ws.Range("M16").Value = -1451.875
$ws.Range("H31").Value = 2183
$ws.Range("I31").Value = 2167.2903
$ws.Range("K31").Value = 2167.2903
$ws.Range("M31").Value = -1872.2903
$ws.Range("H34").Value = 2183
$ws.Range("I34").Value = 2167.2903
$ws.Range("K34").Value = 2167.2903
$ws.Range("M34").Value = -1965.2903
$ws.Range("H99").Value = 3350
$ws.Range("I99").Value = 3307.6924
$ws.Range("J99").Value = 3428.5715
$ws.Range("K99").Value = 3307.6924
$ws.Range("L99").Value = 3428.5715
$ws.Range("M99").Value = -1809.6924
$ws.Range("N99").Value = -6424.5715
$ws.Range("H107").Value = 297.30768
$ws.Range("I107").Value = 154.85715
$ws.Range("J107").Value = 377.08
$ws.Range("K107").Value = 154.85715
$ws.Range("L107").Value = 377.08
$ws.Range("M107").Value = 1765.14285
$ws.Range("N107").Value = -4217.08
$ws.Range("H113").Value = 1741.1
$ws.Range("I113").Value = 1738.875
$ws.Range("K113").Value = 1738.875
$ws.Range("M113").Value = 431.125
$ws.Range("H126").Value = 3350
$ws.Range("I126").Value = 3307.6924
$ws.Range("J126").Value = 3428.5715
$ws.Range("K126").Value = 9923.0772
$ws.Range("L126").Value = 10285.7145
$ws.Range("M126").Value = -7453.0772
$ws.Range("N126").Value = -15225.7145
$ws.Range("H132").Value = 3859.4285
$ws.Range("I132").Value = 2434.5715
$ws.Range("J132").Value = 5284.2856
$ws.Range("K132").Value = 7303.7145
$ws.Range("L132").Value = 15852.8568
$ws.Range("M132").Value = -4773.7145
$ws.Range("N132").Value = -20912.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 484.14285
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 481.5
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 1444.5
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -3940.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5142.857
$ws.Range("H73").Value = 5142.857
$ws.Range("H80").Value = 3124.842
$ws.Range("I80").Value = 3305.1538
$ws.Range("J80").Value = 2734.1667
$ws.Range("K80").Value = 3305.1538
$ws.Range("L80").Value = 2734.1667
$ws.Range("M80").Value = -2307.1538
$ws.Range("N80").Value = -4730.1667
$ws.Range("H83").Value = 3124.842
$ws.Range("I83").Value = 3305.1538
$ws.Range("J83").Value = 2734.1667
$ws.Range("K83").Value = 16525.769
$ws.Range("L83").Value = 13670.8335
$ws.Range("M83").Value = -11533.769
$ws.Range("N83").Value = -23654.8335
$ws.Range("H113").Value = 1900
$ws.Range("I113").Value = 1666.6666
$ws.Range("J113").Value = 2250
$ws.Range("K113").Value = 1666.6666
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 503.3334
$ws.Range("N113").Value = -6590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2421.8667
$ws.Range("I82").Value = 2677.2222
$ws.Range("J82").Value = 2038.8334
$ws.Range("K82").Value = 2677.2222
$ws.Range("L82").Value = 2038.8334
$ws.Range("M82").Value = -2316.2222
$ws.Range("N82").Value = -2760.8334
$ws.Range("H85").Value = 2421.8667
$ws.Range("I85").Value = 2677.2222
$ws.Range("J85").Value = 2038.8334
$ws.Range("K85").Value = 2677.2222
$ws.Range("L85").Value = 2038.8334
$ws.Range("M85").Value = -1429.2222
$ws.Range("N85").Value = -4534.8334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2327.2727
$ws.Range("I81").Value = 2314.2856
$ws.Range("J81").Value = 2350
$ws.Range("K81").Value = 4628.5712
$ws.Range("L81").Value = 4700
$ws.Range("M81").Value = -3567.5712
$ws.Range("N81").Value = -6822
$ws.Range("H84").Value = 2327.2727
$ws.Range("I84").Value = 2314.2856
$ws.Range("J84").Value = 2350
$ws.Range("K84").Value = 23142.856
$ws.Range("L84").Value = 23500
$ws.Range("M84").Value = -17838.856
$ws.Range("N84").Value = -34108
